$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.059.47"
$ws.Range("E2").Value = "  +0.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.264.68"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.31"
$ws.Range("E5").Value = "  +0.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.71"
$ws.Range("E6").Value = "  +2.98%  "

$ws.Range("E7").Value = "  -0.39%  "

$ws.Range("E8").Value = "  -0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  +0.82%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.13"
$ws.Range("E10").Value = "  +7.56%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0790"
$ws.Range("E11").Value = "  -1.02%  "

$ws.Range("E12").Value = "  -0.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.63"
$ws.Range("E13").Value = "  -0.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.611.12"
$ws.Range("E14").Value = "  -0.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.38"
$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.272.59"
$ws.Range("E16").Value = "  -0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.792"
$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.942.00"
$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.41"
$ws.Range("E19").Value = "  -4.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0903"
$ws.Range("E20").Value = "  -0.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.97"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.65"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.81"
$ws.Range("E23").Value = "  -2.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.58"
$ws.Range("E24").Value = "  -0.33%  "

$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.93"
$ws.Range("E26").Value = "  -0.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.68"
$ws.Range("E27").Value = "  -1.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.71"
$ws.Range("E28").Value = "  +5.43%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.50"
$ws.Range("E29").Value = "  -1.23%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.11"
$ws.Range("E30").Value = "  +1.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.22"
$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.23"
$ws.Range("E32").Value = "  -1.86%  "

$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.18"
$ws.Range("E34").Value = "  +4.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0738"
$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.03"
$ws.Range("E36").Value = "  +0.23%  "

$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("E38").Value = "  -1.18%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.83"
$ws.Range("E39").Value = "  +1.91%  "

$ws.Range("E40").Value = "  -1.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.02"
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.41"
$ws.Range("E42").Value = "  +7.01%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.982.52"
$ws.Range("E43").Value = "  -1.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.08"
$ws.Range("E44").Value = "  -3.21%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0283"
$ws.Range("E45").Value = "  +0.25%  "

$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.92"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.92"
$ws.Range("E47").Value = "  -3.55%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.16"
$ws.Range("E48").Value = "  -0.59%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.31"
$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.51"
$ws.Range("E50").Value = "  +0.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "90.93"
$ws.Range("E51").Value = "  -1.04%  "
